# Apply the wording fixes described by the commit diff:
#  - Slide 8: "Extremely lose" -> "Extremely loose"
#  - Slide 9: merge the "Some impedence ... " runs into a single corrected
#             "Some impedance ..." run, and fix "can just errors" -> "can error"

$p = $ppt.ActivePresentation

# --- Slide 8: "Minutae" slide -----------------------------------------
$slide8 = $p.Slides.Item(8)
$body8 = $slide8.Shapes.Item(2).TextFrame.TextRange
$para8_2 = $body8.Paragraphs(2)
$para8_2.Runs(1).Text = "Extremely loose forward-compatible matching is demonstrated (ignore underscore, medial hyphen, spaces, etc..)"

# --- Slide 9: "Matching for C++" slide ---------------------------------
$slide9 = $p.Slides.Item(9)
$body9 = $slide9.Shapes.Item(2).TextFrame.TextRange

# Paragraph 3 currently has 3 runs ("Some " / "impedence" / " mismatch ...").
# Re-set the text twice: the first pass to an unrelated placeholder clears
# the old multi-run split, the second pass assigns the real corrected text
# as a single fresh run (matching the target single-run paragraph).
$para9_3 = $body9.Paragraphs(3)
$para9_3.Text = "."
$para9_3.Text = "Some impedance mismatch with {ZERO-WIDTH-SPACE} versus {ZERO-WIDTH SPACE}"

# Paragraph 4: "Compiler can just errors on inability..." -> "Compiler can error on inability..."
$para9_4 = $body9.Paragraphs(4)
$para9_4.Text = "."
$para9_4.Text = "Compiler can error on inability to find a proper escape sequence to prevent silent spelling errors"
